$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows above the current row 10 ("NoRotation-tilt60deg") ---
# This makes room for the Gaussian-Quadrature scheme (moved up next to the other
# "single orientation" schemes) plus three brand-new spiral schemes.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Data for the new rows: [A value, B label]
$newSchemes = @(
    @(8,  "Gaussian-Quadrature"),
    @(9,  "Spiral-90deg-10rot-5space"),
    @(10, "Spiral-90deg-15rot-5space"),
    @(11, "Spiral-90deg-10rot-3space")
)

$targetRow = 10
foreach ($scheme in $newSchemes) {
    $ws.Cells.Item($targetRow, 1).Value = $scheme[0]
    $ws.Cells.Item($targetRow, 2).Value = $scheme[1]

    # Match the look of column A on the other category rows: bold, centered,
    # top-aligned, thin box border.
    $aCell = $ws.Cells.Item($targetRow, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # Averaged-intensity values of 1 across all HKL columns (C:P)
    for ($c = 3; $c -le 16; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = 1
    }

    $targetRow++
}

# --- Renumber the rows that got pushed down by the insert ---
# They used to be rows 10-15 (A = 8..13); they're now rows 14-19 and need
# column A bumped by 4 (8->12 ... 13->17) to stay sequential.
for ($i = 0; $i -le 5; $i++) {
    $ws.Cells.Item(14 + $i, 1).Value = 12 + $i
}

# --- The old "Gaussian-Quadrature" row (previously last, now pushed to row 20)
# is now a duplicate of the row inserted above, so remove it. ---
$ws.Rows.Item(20).Delete()
